$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

# Insert a new column before column EE, shifting EE:FI -> EF:FJ
$ws.Columns("EE:EE").Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftToRight)

# Populate the newly inserted column EE
$ws.Range("EE1").Value = "02-dec"
$ws.Range("EE2:EE25").Value = "-"

Write-Host "done"
